$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '24.681.19'
$ws.Range("E2").Value = '  +0.11%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.702.06'
$ws.Range("E3").Value = '  +0.38%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9964'
$ws.Range("E4").Value = '  -0.94%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '315.18'
$ws.Range("E5").Value = '  -0.79%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9961'
$ws.Range("E6").Value = '  -0.80%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3977'
$ws.Range("E7").Value = '  +0.13%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4065'
$ws.Range("E8").Value = '  +1.03%  '

$ws.Range("B9").Value = 'Polygon'
$ws.Range("C9").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.520'
$ws.Range("E9").Value = '  +6.42%  '

$ws.Range("B10").Value = 'BinanceUSD'
$ws.Range("C10").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9966'
$ws.Range("E10").Value = '  -1.00%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.32'
$ws.Range("E11").Value = '  +9.00%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08786'
$ws.Range("E12").Value = '  -0.55%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.337'
$ws.Range("E13").Value = '  +9.82%  '

$ws.Range("E14").Value = '  -0.52%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001323'
$ws.Range("E15").Value = '  -0.58%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '7.508'
$ws.Range("E16").Value = '  +3.14%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.701.05'
$ws.Range("E17").Value = '  +0.14%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '101.14'
$ws.Range("E18").Value = '  -1.52%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.07107'
$ws.Range("E19").Value = '  +3.46%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.49'
$ws.Range("E20").Value = '  -1.58%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.751'
$ws.Range("E21").Value = '  -1.54%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9961'
$ws.Range("E22").Value = '  -0.85%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '14.20'
$ws.Range("E23").Value = '  +0.62%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '24.669.78'
$ws.Range("E24").Value = '  +0.12%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.990'
$ws.Range("E25").Value = '  +4.24%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.307'
$ws.Range("E26").Value = '  -0.48%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.42'
$ws.Range("E27").Value = '  +0.32%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '159.09'
$ws.Range("E28").Value = '  -0.36%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.134'
$ws.Range("E29").Value = '  -3.11%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.34'
$ws.Range("E30").Value = '  -0.41%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.475'
$ws.Range("E31").Value = '  +24.28%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.886.88'
$ws.Range("E32").Value = '  +0.33%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.093'
$ws.Range("E33").Value = '  -8.53%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08682'
$ws.Range("E34").Value = '  -2.84%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '7.390'
$ws.Range("E35").Value = '  +19.56%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '11.13'
$ws.Range("E36").Value = '  +0.63%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.946'
$ws.Range("E37").Value = '  +3.66%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2728'
$ws.Range("E38").Value = '  +0.01%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '14.82'
$ws.Range("E39").Value = '  -4.52%  '

$ws.Range("E40").Value = '  +8.10%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.08997'
$ws.Range("E41").Value = '  +0.55%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.475'
$ws.Range("E42").Value = '  +0.70%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7656'
$ws.Range("E43").Value = '  -0.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.7239'
$ws.Range("E44").Value = '  +0.28%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '15.54'
$ws.Range("E45").Value = '  +1.22%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.459'
$ws.Range("E46").Value = '  -1.24%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.173'
$ws.Range("E47").Value = '  +0.47%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.9957'
$ws.Range("E48").Value = '  -0.84%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.51'
$ws.Range("E49").Value = '  -1.13%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.322'
$ws.Range("E50").Value = '  +13.04%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.08019'
$ws.Range("E51").Value = '  +1.00%  '
